# Add two new bash_lib entries to sheet1 ("工作表1"): wc / word count, and sha*sum / shasum x sha1sum.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: wc ---
$ws.Range("A42").Value = "wc"
$ws.Range("B42").Value = "word count"
$ws.Range("C42").Value = "wc is a pipe command, its used to count file's word/line/byte count:`n$ cat file.txt | wc -l  //count by lines`n$ cat file.txt | wc -m //count by char`n$ cat file.txt | wc -c //count by bytes"

$ws.Rows.Item(42).RowHeight = 78.75
$ws.Range("C42").WrapText = $true

# --- Row 43: sha*sum ---
$ws.Range("A43").Value = "sha*sum"
$ws.Range("B43").Value = "shasum x sha1sum"
$ws.Range("C43").Value = "Command shasum is the extense version of sha1sum. Use below shasum syntax to perform sha1sum:`n$ shasum -a 1 {file}`n$ shasum -a 512 {file}   // perform SHA512 algorithm"

$ws.Rows.Item(43).RowHeight = 63
$ws.Range("C43").WrapText = $true

# Update the active selection to reflect the next empty row (matches the
# author's workbook state after typing the new entries: C44).
[void]$ws.Range("C44").Select()
